# Resolved the Date Bug
# Append the two missing time-series rows (52 and 53) that were dropped
# because of the date bug, restoring the Sensor-1 / Sensor-2 readings
# for the two extra timestamps recorded just after midnight and in the
# afternoon of 19/2/2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(52, 1).Value = "19/2/2024, 12:44:37 am IST"
$ws.Cells.Item(52, 2).Value = 270.68928
$ws.Cells.Item(52, 3).Value = 270.49391

$ws.Cells.Item(53, 1).Value = "19/2/2024, 2:02:09 pm IST"
$ws.Cells.Item(53, 2).Value = 116.83392
$ws.Cells.Item(53, 3).Value = 115.39074
